# sds2011.xlsx update:
#  - Added converted GIS data for Scotland
#  - Added 2010 & 2011 mid-year population estimates
#  - Added dzone & izone shape files (converted from ngrid to degrees)
#
# Concretely this translates into the following workbook edits:
#  1. Rename the "summarymatrix" sheet to "summary".
#  2. Add a Scotland-wide totals row (row 16) to the "prevelence" sheet.
#  3. Populate the previously-empty per-health-board statistic rows on the
#     "summarymatrix"/"summary" sheet (rows 4, 5, 6, 31, 32, 35, 36).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename "summarymatrix" -> "summary"
# ------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("summarymatrix")
$wsSummary.Name = "summary"

# ------------------------------------------------------------------
# 2. "prevelence" sheet: add Scotland-wide total row 16
# ------------------------------------------------------------------
$wsPrev = $wb.Worksheets.Item("prevelence")

$wsPrev.Range("C16").Value = "Scotland"

$wsPrev.Range("D16").Formula = "=SUM(D2:D15)"
$wsPrev.Range("E16").Formula = "=SUM(E2:E15)"
$wsPrev.Range("F16").Formula = "=SUM(F2:F15)"
$wsPrev.Range("G16").Formula = "=SUM(G2:G15)"
$wsPrev.Range("H16").Formula = "=SUM(H2:H15)"

$wsPrev.Range("A16:J16").Borders.Item(8).LineStyle = 1
$wsPrev.Range("A16:J16").Borders.Item(9).LineStyle = 1
$wsPrev.Range("D16:H16").NumberFormat = "0"

# ------------------------------------------------------------------
# 3. "summary" sheet: fill in the newly-converted figures
# ------------------------------------------------------------------
$wsSum = $wb.Worksheets.Item("summary")

# Row 4 - Crude prevalence of T1 diabetes
$wsSum.Range("C4").Value = 0.6
$wsSum.Range("D4").Value = 0.5
$wsSum.Range("E4").Value = 0.6
$wsSum.Range("F4").Value = 0.5
$wsSum.Range("G4").Value = 0.6
$wsSum.Range("H4").Value = 0.5
$wsSum.Range("I4").Value = 0.6
$wsSum.Range("J4").Value = 0.6
$wsSum.Range("K4").Value = 0.6
$wsSum.Range("L4").Value = 0.5
$wsSum.Range("M4").Value = 0.6
$wsSum.Range("N4").Value = 0.6
$wsSum.Range("O4").Value = 0.5
$wsSum.Range("P4").Value = 0.7
$wsSum.Range("Q4").Value = 0.5
$wsSum.Range("R4").Value = 0.7
$wsSum.Range("S4").Value = 0.6
$wsSum.Range("T4").Value = 0.6
$wsSum.Range("U4").Value = 0.5
$wsSum.Range("V4").Value = 0.5
$wsSum.Range("W4").Value = 0.6
$wsSum.Range("X4").Value = 0.6
$wsSum.Range("Y4").Value = 0.7

# Row 5 - Crude prevalence of T2 diabetes
$wsSum.Range("C5").Value = 4.9
$wsSum.Range("D5").Value = 4.3
$wsSum.Range("E5").Value = 4.9
$wsSum.Range("F5").Value = 4.4
$wsSum.Range("G5").Value = 4.3
$wsSum.Range("H5").Value = 4.2
$wsSum.Range("I5").Value = 3.8
$wsSum.Range("J5").Value = 4
$wsSum.Range("K5").Value = 4.4
$wsSum.Range("L5").Value = 3.5
$wsSum.Range("M5").Value = 4.2
$wsSum.Range("N5").Value = 3.9
$wsSum.Range("O5").Value = 4.5
$wsSum.Range("P5").Value = 4
$wsSum.Range("Q5").Value = 3.5
$wsSum.Range("R5").Value = 4.9
$wsSum.Range("S5").Value = 4.2
$wsSum.Range("T5").Value = 4.3
$wsSum.Range("U5").Value = 3.5
$wsSum.Range("V5").Value = 4
$wsSum.Range("W5").Value = 4.3
$wsSum.Range("X5").Value = 4.4
$wsSum.Range("Y5").Value = 4.9

# Row 6 - Estimated % pop undiagnosed diabetes
$wsSum.Range("C6").Value = 0.9
$wsSum.Range("D6").Value = 1
$wsSum.Range("E6").Value = 1.6
$wsSum.Range("F6").Value = 0.7
$wsSum.Range("G6").Value = 0.7
$wsSum.Range("H6").Value = 0.8
$wsSum.Range("I6").Value = 0.9
$wsSum.Range("J6").Value = 1.9
$wsSum.Range("K6").Value = 0.8
$wsSum.Range("L6").Value = 0.9
$wsSum.Range("M6").Value = 1.3
$wsSum.Range("N6").Value = 1.9
$wsSum.Range("O6").Value = 0.9
$wsSum.Range("P6").Value = 2.7
$wsSum.Range("Q6").Value = 0.7
$wsSum.Range("R6").Value = 2.7
$wsSum.Range("S6").Value = 1.2
$wsSum.Range("T6").Value = 0.9
$wsSum.Range("U6").Value = 0.7
$wsSum.Range("V6").Value = 0.8
$wsSum.Range("W6").Value = 0.9
$wsSum.Range("X6").Value = 1.5
$wsSum.Range("Y6").Value = 2.7

# Row 31 - % of diabetes pop (T1) obese (BMI >= 30)
$wsSum.Range("C31").Value = 25.4
$wsSum.Range("D31").Value = 29.6
$wsSum.Range("E31").Value = 27.7
$wsSum.Range("F31").Value = 25.9
$wsSum.Range("G31").Value = 23.4
$wsSum.Range("H31").Value = 22.5
$wsSum.Range("I31").Value = 23.1
$wsSum.Range("J31").Value = 25.1
$wsSum.Range("K31").Value = 26.4
$wsSum.Range("L31").Value = 23.2
$wsSum.Range("M31").Value = 29.1
$wsSum.Range("N31").Value = 32.1
$wsSum.Range("O31").Value = 23.3
$wsSum.Range("P31").Value = 24.2
$wsSum.Range("Q31").Value = 22.5
$wsSum.Range("R31").Value = 32.1
$wsSum.Range("S31").Value = 25.8
$wsSum.Range("T31").Value = 25.3
$wsSum.Range("U31").Value = 22.5
$wsSum.Range("V31").Value = 23.3
$wsSum.Range("W31").Value = 25.3
$wsSum.Range("X31").Value = 27.4
$wsSum.Range("Y31").Value = 32.1

# Row 32 - % of diabetes pop (T2) obese (BMI >= 30)
$wsSum.Range("C32").Value = 54.8
$wsSum.Range("D32").Value = 55.6
$wsSum.Range("E32").Value = 55
$wsSum.Range("F32").Value = 58.7
$wsSum.Range("G32").Value = 57.4
$wsSum.Range("H32").Value = 54.7
$wsSum.Range("I32").Value = 53.1
$wsSum.Range("J32").Value = 55.2
$wsSum.Range("K32").Value = 56.7
$wsSum.Range("L32").Value = 56
$wsSum.Range("M32").Value = 60.3
$wsSum.Range("N32").Value = 60.9
$wsSum.Range("O32").Value = 55.5
$wsSum.Range("P32").Value = 57.5
$wsSum.Range("Q32").Value = 53.1
$wsSum.Range("R32").Value = 60.9
$wsSum.Range("S32").Value = 56.5
$wsSum.Range("T32").Value = 55.8
$wsSum.Range("U32").Value = 53.1
$wsSum.Range("V32").Value = 55.1
$wsSum.Range("W32").Value = 55.8
$wsSum.Range("X32").Value = 57.5
$wsSum.Range("Y32").Value = 60.9

# Row 35 - % of diabetes pop (T1) with HbA1c > 75mmol/mol
$wsSum.Range("C35").Value = 37.6
$wsSum.Range("D35").Value = 37.2
$wsSum.Range("E35").Value = 32
$wsSum.Range("F35").Value = 37.9
$wsSum.Range("G35").Value = 40.9
$wsSum.Range("H35").Value = 46.9
$wsSum.Range("I35").Value = 35.6
$wsSum.Range("J35").Value = 36.4
$wsSum.Range("K35").Value = 34.7
$wsSum.Range("L35").Value = 31.1
$wsSum.Range("M35").Value = 33.6
$wsSum.Range("N35").Value = 26.7
$wsSum.Range("O35").Value = 43.6
$wsSum.Range("P35").Value = 40
$wsSum.Range("Q35").Value = 26.7
$wsSum.Range("R35").Value = 46.9
$wsSum.Range("S35").Value = 36.7
$wsSum.Range("T35").Value = 36.8
$wsSum.Range("U35").Value = 26.7
$wsSum.Range("V35").Value = 33.9
$wsSum.Range("W35").Value = 36.8
$wsSum.Range("X35").Value = 39.5
$wsSum.Range("Y35").Value = 46.9

# Row 36 - % of diabetes pop (T2) with HbA1c > 75mmol/mol
$wsSum.Range("C36").Value = 13.4
$wsSum.Range("D36").Value = 11.5
$wsSum.Range("E36").Value = 12.2
$wsSum.Range("F36").Value = 13
$wsSum.Range("G36").Value = 14.3
$wsSum.Range("H36").Value = 16.6
$wsSum.Range("I36").Value = 16.2
$wsSum.Range("J36").Value = 14.9
$wsSum.Range("K36").Value = 14.7
$wsSum.Range("L36").Value = 12.7
$wsSum.Range("M36").Value = 16.9
$wsSum.Range("N36").Value = 14.4
$wsSum.Range("O36").Value = 14.8
$wsSum.Range("P36").Value = 18.3
$wsSum.Range("Q36").Value = 11.5
$wsSum.Range("R36").Value = 18.3
$wsSum.Range("S36").Value = 14.6
$wsSum.Range("T36").Value = 14.6
$wsSum.Range("U36").Value = 11.5
$wsSum.Range("V36").Value = 13.1
$wsSum.Range("W36").Value = 14.6
$wsSum.Range("X36").Value = 15.9
$wsSum.Range("Y36").Value = 18.3

# ------------------------------------------------------------------
# Restore cell selections on every sheet, then leave "summary" active,
# matching the saved file.
# ------------------------------------------------------------------
$wsPrev.Select()
$wsPrev.Range("I18").Select()

$wsFact = $wb.Worksheets.Item("fact")
$wsFact.Select()
$wsFact.Range("B30").Select()

$wsSum.Select()
$wsSum.Range("B31").Select()
